# Auto-generated edit script: update Leve profit calculations across sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 58824140
$ws.Range("I98").Value = 58824140
$ws.Range("K98").Value = 58824140
$ws.Range("M98").Value = -58822642

$ws.Range("H122").Value = 58824140
$ws.Range("I122").Value = 58824140
$ws.Range("K122").Value = 176472420
$ws.Range("M122").Value = -176469970

$ws.Range("H131").Value = 1725.7778
$ws.Range("I131").Value = 1741.625
$ws.Range("K131").Value = 5224.875
$ws.Range("M131").Value = -184.875

$ws.Range("H132").Value = 2734.7058
$ws.Range("I132").Value = 2726.6667
$ws.Range("K132").Value = 8180.000100000001
$ws.Range("M132").Value = -5650.000100000001

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").Value = $null

$ws.Range("H135").Value = 15101.143
$ws.Range("I135").Value = 902.6667
$ws.Range("K135").Value = 8124.0003
$ws.Range("M135").Value = -5589.0003

$ws.Range("H137").Value = 3236.0312
$ws.Range("I137").Value = 1256.52
$ws.Range("J137").Value = 10305.714
$ws.Range("K137").Value = 3769.56
$ws.Range("L137").Value = 30917.142
$ws.Range("M137").Value = -1219.56
$ws.Range("N137").Value = -36017.142

$ws.Range("H139").Value = 88000
$ws.Range("J139").Value = 88000
$ws.Range("L139").Value = 88000
$ws.Range("N139").Value = -98280

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9805102
$ws.Range("I32").Value = 9805102
$ws.Range("K32").Value = 9805102
$ws.Range("M32").Value = -9804815

$ws.Range("H45").Value = 2750.7144
$ws.Range("I45").Value = 2654.6155
$ws.Range("K45").Value = 2654.6155
$ws.Range("M45").Value = -2277.6155

$ws.Range("H61").Value = 23863644
$ws.Range("I61").Value = 33338038
$ws.Range("J61").Value = 177658.5
$ws.Range("K61").Value = 33338038
$ws.Range("L61").Value = 177658.5
$ws.Range("M61").Value = -33337826
$ws.Range("N61").Value = -178082.5

$ws.Range("H74").Value = 5820020
$ws.Range("I74").Value = 7814082.5
$ws.Range("K74").Value = 7814082.5
$ws.Range("M74").Value = -7813208.5

$ws.Range("H77").Value = 5820020
$ws.Range("I77").Value = 7814082.5
$ws.Range("K77").Value = 39070412.5
$ws.Range("M77").Value = -39066044.5

$ws.Range("H97").Value = 1266.3684
$ws.Range("I97").Value = 1210.125
$ws.Range("J97").Value = 1566.3334
$ws.Range("K97").Value = 1210.125
$ws.Range("L97").Value = 1566.3334
$ws.Range("M97").Value = -714.125
$ws.Range("N97").Value = -2558.3334

$ws.Range("H110").Value = 2411.4546
$ws.Range("I110").Value = 2250.4211
$ws.Range("K110").Value = 2250.4211
$ws.Range("M110").Value = -205.4211

$ws.Range("H112").Value = 23457.4
$ws.Range("J112").Value = 23457.4
$ws.Range("L112").Value = 23457.4
$ws.Range("N112").Value = -26411.4

$ws.Range("H122").Value = 1226.6364
$ws.Range("I122").Value = 1249.3
$ws.Range("K122").Value = 3747.9
$ws.Range("M122").Value = -1297.9

$ws.Range("H132").Value = 4723.8687
$ws.Range("I132").Value = 3106.4915
$ws.Range("J132").Value = 10337.117
$ws.Range("K132").Value = 9319.4745
$ws.Range("L132").Value = 31011.351
$ws.Range("M132").Value = -6789.4745
$ws.Range("N132").Value = -36071.351

$ws.Range("H136").Value = 23863644
$ws.Range("I136").Value = 33338038
$ws.Range("J136").Value = 177658.5
$ws.Range("K136").Value = 100014114
$ws.Range("L136").Value = 532975.5
$ws.Range("M136").Value = -100011564
$ws.Range("N136").Value = -538075.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4571.6
$ws.Range("I20").Value = 3857.4443
$ws.Range("J20").Value = 10999
$ws.Range("K20").Value = 3857.4443
$ws.Range("L20").Value = 10999
$ws.Range("M20").Value = -3610.4443
$ws.Range("N20").Value = -11493

$ws.Range("H94").Value = 1490.2667
$ws.Range("I94").Value = 1613.3334
$ws.Range("J94").Value = 998
$ws.Range("K94").Value = 1613.3334
$ws.Range("L94").Value = 998
$ws.Range("M94").Value = -1162.3334
$ws.Range("N94").Value = -1900

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 2053.2942
$ws.Range("I7").Value = 141.64285
$ws.Range("J7").Value = 10974.333
$ws.Range("K7").Value = 141.64285
$ws.Range("L7").Value = 10974.333
$ws.Range("M7").Value = -28.64285000000001
$ws.Range("N7").Value = -11200.333

$ws.Range("H81").Value = 123000
$ws.Range("J81").Value = 123000
$ws.Range("L81").Value = 123000
$ws.Range("N81").Value = -124996

$ws.Range("H84").Value = 123000
$ws.Range("J84").Value = 123000
$ws.Range("L84").Value = 369000
$ws.Range("N84").Value = -378984

$ws.Range("H132").Value = 3959
$ws.Range("I132").Value = 4112.857
$ws.Range("J132").Value = 3600
$ws.Range("K132").Value = 12338.571
$ws.Range("L132").Value = 10800
$ws.Range("M132").Value = -9808.571
$ws.Range("N132").Value = -15860

$ws.Range("H134").Value = 288994.88
$ws.Range("I134").Value = 371953.28
$ws.Range("K134").Value = 1115859.84
$ws.Range("M134").Value = -1113324.84

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 1184.6923
$ws.Range("I33").Value = 1227.3636
$ws.Range("K33").Value = 7364.1816
$ws.Range("M33").Value = -7081.1816

$ws.Range("H125").Value = 13433.286
$ws.Range("J125").Value = 13433.286
$ws.Range("L125").Value = 40299.858
$ws.Range("N125").Value = -50139.858

$ws.Range("H132").Value = 1618.7778
$ws.Range("I132").Value = 1175
$ws.Range("J132").Value = 2173.5
$ws.Range("K132").Value = 10575
$ws.Range("L132").Value = 19561.5
$ws.Range("M132").Value = -8045
$ws.Range("N132").Value = -24621.5

$ws.Range("H133").Value = 5104.4
$ws.Range("I133").Value = 4893.778
$ws.Range("K133").Value = 14681.334
$ws.Range("M133").Value = -9621.334000000001

$ws.Range("H134").Value = 3401.4333
$ws.Range("I134").Value = 2091.8696
$ws.Range("J134").Value = 7704.2856
$ws.Range("K134").Value = 6275.6088
$ws.Range("L134").Value = 23112.8568
$ws.Range("M134").Value = -1205.6088
$ws.Range("N134").Value = -33252.8568

$ws.Range("H139").Value = 2305.5264

$ws.Range("H140").Value = 337173.78
$ws.Range("I140").Value = 432347.84
$ws.Range("K140").Value = 1297043.52
$ws.Range("M140").Value = -1291863.52

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1612.9375
$ws.Range("I107").Value = 1127
$ws.Range("J107").Value = 2682
$ws.Range("K107").Value = 1127
$ws.Range("L107").Value = 2682
$ws.Range("M107").Value = 793
$ws.Range("N107").Value = -6522

$ws.Range("H110").Value = 129000
$ws.Range("J110").Value = 129000
$ws.Range("L110").Value = 129000
$ws.Range("N110").Value = -137180

$ws.Range("H111").Value = 52409.4
$ws.Range("J111").Value = 52409.4
$ws.Range("L111").Value = 52409.4
$ws.Range("N111").Value = -58543.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2883.9285
$ws.Range("I16").Value = 2609.5557
$ws.Range("J16").Value = 3377.8
$ws.Range("K16").Value = 2609.5557
$ws.Range("L16").Value = 3377.8
$ws.Range("M16").Value = -2439.5557
$ws.Range("N16").Value = -3717.8

$ws.Range("H22").Value = 3749
$ws.Range("I22").Value = 5301.5
$ws.Range("K22").Value = 5301.5
$ws.Range("M22").Value = -5006.5

$ws.Range("H27").Value = 3749
$ws.Range("I27").Value = 5301.5
$ws.Range("K27").Value = 5301.5
$ws.Range("M27").Value = -5194.5

$ws.Range("H61").Value = 726.8570999999999
$ws.Range("J61").Value = 998.5
$ws.Range("L61").Value = 998.5
$ws.Range("N61").Value = -1402.5

$ws.Range("H68").Value = 1399.6
$ws.Range("I68").Value = 1499.75
$ws.Range("K68").Value = 1499.75
$ws.Range("M68").Value = -750.75

$ws.Range("H71").Value = 1399.6
$ws.Range("I71").Value = 1499.75
$ws.Range("K71").Value = 7498.75
$ws.Range("M71").Value = -3754.75

$ws.Range("H110").Value = 76711.336
$ws.Range("J110").Value = 76711.336
$ws.Range("L110").Value = 76711.336
$ws.Range("N110").Value = -84891.336

$ws.Range("H113").Value = 726.8570999999999
$ws.Range("J113").Value = 998.5
$ws.Range("L113").Value = 998.5
$ws.Range("N113").Value = -5338.5

$ws.Range("H136").Value = 48952.883
$ws.Range("I136").Value = 5224
$ws.Range("K136").Value = 15672
$ws.Range("M136").Value = -13122

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 10000100
$ws.Range("J5").Value = 10000100
$ws.Range("L5").Value = 10000100
$ws.Range("N5").Value = -10000324

$ws.Range("H11").Value = 38012
$ws.Range("J11").Value = 38012
$ws.Range("L11").Value = 38012
$ws.Range("N11").Value = -38296

$ws.Range("H132").Value = 13954.333
$ws.Range("J132").Value = 29234.309
$ws.Range("L132").Value = 87702.927
$ws.Range("N132").Value = -92762.927

$ws.Range("H133").Value = 89807.5
$ws.Range("J133").Value = 89807.5
$ws.Range("L133").Value = 89807.5
$ws.Range("N133").Value = -99927.5

$ws.Range("H136").Value = 12370.083
$ws.Range("I136").Value = 1493.7858
$ws.Range("K136").Value = 4481.357400000001
$ws.Range("M136").Value = -1931.357400000001

$ws.Range("H137").Value = 105998.336
$ws.Range("J137").Value = 105998.336
$ws.Range("L137").Value = 105998.336
$ws.Range("N137").Value = -116198.336
